# Fix and sort Excel files
# Insert a new county row for "Anchorage" / "AK_Anchorage" into the
# alphabetically-sorted list of Alaska counties (it belongs right after
# "Aleutians West" in row 3 and before "Bethel", which currently sits on
# row 4), pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 currently holds "Bethel" / "AK_Bethel" - push it (and everything
# below it) down one row to make room for the new entry.
$ws.Rows("4:4").Insert()

# Populate the newly inserted row with the Anchorage data.
$ws.Cells.Item(4, 1).Value2 = "Anchorage"
$ws.Cells.Item(4, 2).Value2 = "AK_Anchorage"
